$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns for the LowRate2 scenario
$ws.Range("H1").Value = "MA.TDA_MA_multiTier_TDAamortAS_OYLM_TDA_LowRate2_base"
$ws.Range("I1").Value = "MA.TDA_MA_multiTier_TDAamortAS_OYLM_TDA_LowRate2_lowG"

# New data columns matching the 3 existing data rows
$ws.Range("H2").Value = 0.4363754059958826
$ws.Range("I2").Value = 0.4363754059958826

$ws.Range("H3").Value = 0.22876177519195645
$ws.Range("I3").Value = 0.19612304652976992

$ws.Range("H4").Value = 0.18990601239107158
$ws.Range("I4").Value = 0.12535359893479336
